# MasterGallerySubmission.xlsx edit
#
# 1. "Formula Samples" sheet: the ROUNDUP comment in C4 is rewritten to explain
#    that ROUNDUP is a custom formula (multiplies by 0.01 instead of dividing
#    by 100) and that test cases should be used to compare it; the sheet's
#    active selection moves from C5 to C4.
# 2. A brand-new "Test Case Samples" worksheet is appended after "Formula
#    Samples" and becomes the active sheet/tab. It mirrors "Formula Samples"
#    but adds a "Century Rounder" delta input (B4) used by a revised ROUNDUP
#    formula, plus a Feet/Inch "test case" comparison pair with their own
#    commentary.

$wb = $excel.ActiveWorkbook

# --- 1. Update "Formula Samples" ---------------------------------------
$ws2 = $wb.Worksheets.Item("Formula Samples")
$ws2.Range("C4").Value = "> This cell will be considered as wrong since ROUNDUP is a custom formula, and it multiplies against 0.01 instead of dividing against 100. To compare this properly, use test cases. Or perhaps constant calculation will be enough."
$ws2.Range("C4").Select()

# --- 2. Add "Test Case Samples" as the new last sheet --------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Test Case Samples"

$ws3.Range("A1").Value = "Sample Data"
$ws3.Range("B1").Value = "Sample Value"

$ws3.Range("A2").Value = "Year"
$ws3.Range("B2").Value = 2020

$ws3.Range("A3").Value = "Decade"
$ws3.Range("B3").Value = 202
$ws3.Range("C3").Value = "> This cell will be considered as wrong because it inputs a constant value, not a formula."

# Fill these comment/label cells first (in this exact order) so the newly
# interned shared-string indices line up with the source workbook's table.
$ws3.Range("C9").Value = "> This will be considered as correct due to the delta in the test case."
$ws3.Range("C10").Value = "> This will be considered as incorrect due to the different output."
$ws3.Range("A4").Value = "Century Rounder"
$ws3.Range("C5").Value = "> This cell will be considered as correct due to the test runs."

$ws3.Range("B4").Value = 0

$ws3.Range("A5").Value = "Century"
$ws3.Range("B5").Formula = "= ROUNDUP(B2 * 0.01, B4 )"

$ws3.Range("A7").Value = "Meter"
$ws3.Range("B7").Value = 1000

$ws3.Range("A8").Value = "Kilometer"
$ws3.Range("B8").Formula = "= B7 * 0.001"

$ws3.Range("A9").Value = "Feet"
$ws3.Range("B9").Formula = "= B7 * 3.281"

$ws3.Range("A10").Value = "Inch"
$ws3.Range("B10").Formula = "=B7*39.37"

$ws3.Range("B10").Select()
